# Add three new columns (AD: Wins, AE: Losses, AF: Ties) holding the team's
# season record, which was missing from the original scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - same header style as the rest of row 1 (bold, bordered,
# centered/top aligned). Copy the formatting from an existing header cell
# so we reuse the workbook's existing style record instead of minting a
# new one.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-51: every player on the roster shares the team's 58-104-0
# season record.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 58
    $ws.Cells.Item($r, 31).Value = 104
    $ws.Cells.Item($r, 32).Value = 0
}
